$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24 (shifts Joan Garcia -> 25, Szczesny -> 26)
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with Jofre Torrents' betting/analytics stats
$ws.Cells.Item(24, 1).Value = 'La_Liga'
$ws.Cells.Item(24, 2).Value = 'Barcelona'
$ws.Cells.Item(24, 3).Value = 'Jofre Torrents '
$ws.Cells.Item(24, 4).Value = 1587196
$ws.Cells.Item(24, 5).Value = 6.4
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 47
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0
$ws.Cells.Item(24, 18).Value = 0
$ws.Cells.Item(24, 19).Value = 0
$ws.Cells.Item(24, 20).Value = 0
$ws.Cells.Item(24, 21).Value = 0
$ws.Cells.Item(24, 22).Value = 0
$ws.Cells.Item(24, 23).Value = 0
$ws.Cells.Item(24, 24).Value = 0
$ws.Cells.Item(24, 25).Value = 0
$ws.Cells.Item(24, 26).Value = 0
$ws.Cells.Item(24, 27).Value = 0.0353411
$ws.Cells.Item(24, 28).Value = 54
$ws.Cells.Item(24, 29).Value = 0
$ws.Cells.Item(24, 30).Value = 0
$ws.Cells.Item(24, 31).Value = 30
$ws.Cells.Item(24, 32).Value = 88.235294117647
$ws.Cells.Item(24, 33).Value = 34
$ws.Cells.Item(24, 34).Value = 11
$ws.Cells.Item(24, 35).Value = 19
$ws.Cells.Item(24, 36).Value = 15
$ws.Cells.Item(24, 37).Value = 1
$ws.Cells.Item(24, 38).Value = 100
$ws.Cells.Item(24, 39).Value = 1
$ws.Cells.Item(24, 40).Value = 33.333333333333
$ws.Cells.Item(24, 41).Value = 0
$ws.Cells.Item(24, 42).Value = 2
$ws.Cells.Item(24, 43).Value = 0
$ws.Cells.Item(24, 44).Value = 3
$ws.Cells.Item(24, 45).Value = 0
$ws.Cells.Item(24, 46).Value = 2
$ws.Cells.Item(24, 47).Value = 0
$ws.Cells.Item(24, 48).Value = 0
$ws.Cells.Item(24, 49).Value = 0
$ws.Cells.Item(24, 50).Value = 0
$ws.Cells.Item(24, 51).Value = 0
$ws.Cells.Item(24, 52).Value = 0
$ws.Cells.Item(24, 53).Value = 3
$ws.Cells.Item(24, 54).Value = 37.5
$ws.Cells.Item(24, 55).Value = 2
$ws.Cells.Item(24, 56).Value = 50
$ws.Cells.Item(24, 57).Value = 1
$ws.Cells.Item(24, 58).Value = 25
$ws.Cells.Item(24, 59).Value = 10
$ws.Cells.Item(24, 60).Value = 0
$ws.Cells.Item(24, 61).Value = 0
$ws.Cells.Item(24, 62).Value = 0
$ws.Cells.Item(24, 63).Value = 0
$ws.Cells.Item(24, 64).Value = 0
$ws.Cells.Item(24, 65).Value = 0
$ws.Cells.Item(24, 66).Value = 0
$ws.Cells.Item(24, 67).Value = 2
$ws.Cells.Item(24, 68).Value = 0
$ws.Cells.Item(24, 69).Value = 19.2
$ws.Cells.Item(24, 70).Value = 3
$ws.Cells.Item(24, 71).Value = 0
$ws.Cells.Item(24, 72).Value = 4
$ws.Cells.Item(24, 73).Value = 0
$ws.Cells.Item(24, 74).Value = 0
$ws.Cells.Item(24, 75).Value = 0
$ws.Cells.Item(24, 76).Value = 0
$ws.Cells.Item(24, 77).Value = 0
$ws.Cells.Item(24, 78).Value = 0
$ws.Cells.Item(24, 79).Value = 0
$ws.Cells.Item(24, 80).Value = 1
$ws.Cells.Item(24, 81).Value = 0
$ws.Cells.Item(24, 82).Value = 0
$ws.Cells.Item(24, 83).Value = 0
$ws.Cells.Item(24, 84).Value = 0
$ws.Cells.Item(24, 85).Value = 0
$ws.Cells.Item(24, 86).Value = 0
$ws.Cells.Item(24, 87).Value = 0
$ws.Cells.Item(24, 88).Value = 2
$ws.Cells.Item(24, 89).Value = 0
$ws.Cells.Item(24, 90).Value = 0
$ws.Cells.Item(24, 91).Value = 0
$ws.Cells.Item(24, 92).Value = 0
$ws.Cells.Item(24, 93).Value = 0
$ws.Cells.Item(24, 94).Value = 0
$ws.Cells.Item(24, 95).Value = 0
$ws.Cells.Item(24, 96).Value = 0
$ws.Cells.Item(24, 97).Value = 1
$ws.Cells.Item(24, 98).Value = 3
$ws.Cells.Item(24, 99).Value = 5
$ws.Cells.Item(24, 100).Value = 3
$ws.Cells.Item(24, 101).Value = 0
$ws.Cells.Item(24, 102).Value = 0
$ws.Cells.Item(24, 103).Value = 0
$ws.Cells.Item(24, 104).Value = 1
$ws.Cells.Item(24, 105).Value = 2
$ws.Cells.Item(24, 106).Value = 100
$ws.Cells.Item(24, 107).Value = 0
$ws.Cells.Item(24, 108).Value = 0
$ws.Cells.Item(24, 109).Value = 0
$ws.Cells.Item(24, 110).Value = 12
$ws.Cells.Item(24, 111).Value = 22
$ws.Cells.Item(24, 112).Value = 0
$ws.Cells.Item(24, 113).Value = 0
$ws.Cells.Item(24, 114).Value = 2139676
$ws.Cells.Item(24, 115).Value = 'overall'

$wb.Save()
